$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.1832435933191248
$ws.Range("D2").Value = 0.8562851856351483

$ws.Range("C3").Value = 0.6983457689830992
$ws.Range("D3").Value = 0.4922774486340598

$ws.Range("C4").Value = 0.4500564402174905
$ws.Range("D4").Value = 0.6570731942073555

$ws.Range("C5").Value = 2.167558588652516
$ws.Range("D5").Value = 0.0412875051353867

$ws.Range("C6").Value = 0.5126421167307502
$ws.Range("D6").Value = 0.6133106005456153

$ws.Range("C7").Value = 0.4030493446970735
$ws.Range("D7").Value = 0.6908020548691507

$ws.Range("C8").Value = 2.24597874705914
$ws.Range("D8").Value = 0.03507759722231163

$ws.Range("C9").Value = -0.3527997888632249
$ws.Range("D9").Value = 0.727596272958678

$ws.Range("C10").Value = 1.173589285924658
$ws.Range("D10").Value = 0.2531089427007296

$ws.Range("C11").Value = 1.723521432056095
$ws.Range("D11").Value = 0.09882166873345533
$ws.Range("G11").Value = "No"
